{"js": "// Update the answer grid: each cell's \"A\u00d7B=C\" text is replaced with a new\n// equation, per the commit's regenerated answer key. Matches are exact,\n// case-sensitive, whole strings, so a body.search()+insertText() per pair\n// is unambiguous (every \"old\" string is unique in the document).\n\nconst replacements = [\n  [\"35\u00d716=560\", \"83\u00d721=1743\"],\n  [\"65\u00d739=2535\", \"25\u00d763=1575\"],\n  [\"97\u00d737=3589\", \"11\u00d750=550\"],\n  [\"45\u00d784=3780\", \"13\u00d754=702\"],\n  [\"46\u00d766=3036\", \"48\u00d718=864\"],\n  [\"98\u00d776=7448\", \"29\u00d726=754\"],\n  [\"69\u00d718=1242\", \"53\u00d783=4399\"],\n  [\"28\u00d744=1232\", \"52\u00d733=1716\"],\n  [\"38\u00d762=2356\", \"13\u00d752=676\"],\n  [\"33\u00d771=2343\", \"64\u00d793=5952\"],\n  [\"50\u00d775=3750\", \"52\u00d798=5096\"],\n  [\"61\u00d753=3233\", \"87\u00d752=4524\"],\n  [\"80\u00d712=960\", \"33\u00d720=660\"],\n  [\"57\u00d751=2907\", \"55\u00d775=4125\"],\n  [\"16\u00d721=336\", \"61\u00d727=1647\"],\n  [\"35\u00d763=2205\", \"54\u00d736=1944\"],\n  [\"49\u00d723=1127\", \"13\u00d747=611\"],\n  [\"28\u00d713=364\", \"92\u00d762=5704\"],\n  [\"62\u00d747=2914\", \"84\u00d782=6888\"],\n  [\"74\u00d740=2960\", \"26\u00d789=2314\"],\n  [\"63\u00d779=4977\", \"81\u00d744=3564\"],\n  [\"60\u00d792=5520\", \"52\u00d767=3484\"],\n  [\"18\u00d715=270\", \"48\u00d789=4272\"],\n  [\"52\u00d753=2756\", \"92\u00d715=1380\"],\n  [\"53\u00d738=2014\", \"27\u00d772=1944\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the answer grid: each cell's \"A\u00d7B=C\" text is replaced with a new\n# equation, per the commit's regenerated answer key. Matches are exact,\n# case-sensitive, whole strings, so Find/Replace per pair is unambiguous.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"35\u00d716=560\";  New = \"83\u00d721=1743\" },\n    @{ Old = \"65\u00d739=2535\"; New = \"25\u00d763=1575\" },\n    @{ Old = \"97\u00d737=3589\"; New = \"11\u00d750=550\" },\n    @{ Old = \"45\u00d784=3780\"; New = \"13\u00d754=702\" },\n    @{ Old = \"46\u00d766=3036\"; New = \"48\u00d718=864\" },\n    @{ Old = \"98\u00d776=7448\"; New = \"29\u00d726=754\" },\n    @{ Old = \"69\u00d718=1242\"; New = \"53\u00d783=4399\" },\n    @{ Old = \"28\u00d744=1232\"; New = \"52\u00d733=1716\" },\n    @{ Old = \"38\u00d762=2356\"; New = \"13\u00d752=676\" },\n    @{ Old = \"33\u00d771=2343\"; New = \"64\u00d793=5952\" },\n    @{ Old = \"50\u00d775=3750\"; New = \"52\u00d798=5096\" },\n    @{ Old = \"61\u00d753=3233\"; New = \"87\u00d752=4524\" },\n    @{ Old = \"80\u00d712=960\";  New = \"33\u00d720=660\" },\n    @{ Old = \"57\u00d751=2907\"; New = \"55\u00d775=4125\" },\n    @{ Old = \"16\u00d721=336\";  New = \"61\u00d727=1647\" },\n    @{ Old = \"35\u00d763=2205\"; New = \"54\u00d736=1944\" },\n    @{ Old = \"49\u00d723=1127\"; New = \"13\u00d747=611\" },\n    @{ Old = \"28\u00d713=364\";  New = \"92\u00d762=5704\" },\n    @{ Old = \"62\u00d747=2914\"; New = \"84\u00d782=6888\" },\n    @{ Old = \"74\u00d740=2960\"; New = \"26\u00d789=2314\" },\n    @{ Old = \"63\u00d779=4977\"; New = \"81\u00d744=3564\" },\n    @{ Old = \"60\u00d792=5520\"; New = \"52\u00d767=3484\" },\n    @{ Old = \"18\u00d715=270\";  New = \"48\u00d789=4272\" },\n    @{ Old = \"52\u00d753=2756\"; New = \"92\u00d715=1380\" },\n    @{ Old = \"53\u00d738=2014\"; New = \"27\u00d772=1944\" }\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $r.Old\n    $find.Replacement.Text = $r.New\n    $find.Execute([ref]$r.Old, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$r.New, 2) | Out-Null\n}\n"}
